# model-flow.pptx — "Added initial baseline qualitative forecasts (v0.14)"
#
# The "FFR, SOFR Qual Forecasts" label box on slide 2 is re-worded to
# "FFR, SOFR Baseline Forecasts", with the middle word ("Baseline")
# kept as its own run (matching the canonical OOXML, which separates
# "FFR, " / "SOFR Baseline " / "Forecasts" into three <a:r> runs, all
# sz=1200 b=1).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the "Rectangle 35" label shape by name (more robust than a
# hard-coded z-order index).
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Rectangle 35") {
        $target = $cand
        break
    }
}

$tr = $target.TextFrame.TextRange

# Replace the whole string (this keeps the existing run-level
# formatting: lang="en-US" sz="1200" b="1").
$tr.Text = "FFR, SOFR Baseline Forecasts"
$tr.Font.Size = 12
$tr.Font.Bold = $true

# Split out "SOFR Baseline " (characters 6-19) into its own run, as in
# the target markup, re-asserting bold on that middle run.
$mid = $tr.Characters(6, 14)
$mid.Font.Bold = $true

Write-Output "Updated '$($target.Name)' text -> '$($tr.Text)'"
